$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 12; existing rows 12-78 shift down to 13-79
$ws.Rows.Item(12).Insert()

# Populate the newly inserted row 12 with a new weekly record
# (values mirror the prior row 12 entry except Fecha (D) and Volumen (M))
$ws.Cells.Item(12, 1).Value = 3
$ws.Cells.Item(12, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(12, 3).Value = "Coquimbo"
$ws.Cells.Item(12, 4).Value = 45133
$ws.Cells.Item(12, 5).Value = 5
$ws.Cells.Item(12, 6).Value = "Fruta"
$ws.Cells.Item(12, 7).Value = 100108
$ws.Cells.Item(12, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(12, 9).Value = 100108004
$ws.Cells.Item(12, 10).Value = "Papaya"
$ws.Cells.Item(12, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(12, 12).Value = "Primera"
$ws.Cells.Item(12, 13).Value = 56
$ws.Cells.Item(12, 14).Value = 20000
$ws.Cells.Item(12, 15).Value = 20000
$ws.Cells.Item(12, 16).Value = 20000
$ws.Cells.Item(12, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(12, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(12, 19).Value = 2000
$ws.Cells.Item(12, 20).Value = 10
